# HFP_template.xlsx update — "aging" the yearly planning sheets by one year
# (drop the elapsed year 2020, add 2057 fully-formatted, clear the bold-ish
# style that had accidentally been applied to the plain "year" column), and
# inserting a "year" column into the Fixed Assets sheet (matching the Debts
# sheet's layout).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Kim" sheet: years 2020-2057 (rows 2-39) -> years 2021-2057 (rows 2-38)
# ---------------------------------------------------------------------
$kim = $wb.Worksheets.Item("Kim")

# Drop the elapsed year (2020) — everything below shifts up one row.
$kim.Rows.Item(2).Delete()

# After the shift, row 36 (now year 2055) lost its B:I currency-formatted
# cells (they used to belong to the old row 37, which only ever had column
# A filled in). Re-apply the same formatting as the row above it so 2055
# looks like every other "normal" year row.
$kim.Range("B35:I35").Copy()
$kim.Range("B36:I36").PasteSpecial(-4122)
$kim.Application.CutCopyMode = 0

# The plain "year" numbers in column A were styled with the bold header
# font by mistake — reset them back to the workbook's Normal style.
$kim.Range("A2:A38").Style = "Normal"

$kim.Range("B7").Select()

# ---------------------------------------------------------------------
# "Sam" sheet: same re-basing, years 2020-2057 (rows 2-39) -> 2021-2057
# (rows 2-38)
# ---------------------------------------------------------------------
$sam = $wb.Worksheets.Item("Sam")

$sam.Rows.Item(2).Delete()

# Row 38 (now year 2057) inherited the old row 39, which never had its
# B:I cells filled in — copy the formatting from the row above so the
# last year matches the rest of the table.
$sam.Range("B37:I37").Copy()
$sam.Range("B38:I38").PasteSpecial(-4122)
$sam.Application.CutCopyMode = 0

$sam.Range("A2:A38").Style = "Normal"

$sam.Range("B5").Select()

# ---------------------------------------------------------------------
# "Debts" sheet: no data changes, just refresh the selection.
# ---------------------------------------------------------------------
$debts = $wb.Worksheets.Item("Debts")
$debts.Columns.Item(4).Select()

# ---------------------------------------------------------------------
# "Fixed Assets" sheet: insert a "year" column (to match "Debts"),
# shifting basis/value/rate/yod/commission one column to the right.
# ---------------------------------------------------------------------
$fa = $wb.Worksheets.Item("Fixed Assets")
$fa.Columns.Item(4).Insert()
$fa.Cells.Item(1, 4).Value = "year"

$fa.Columns.Item(4).Select()

# Fixed Assets is the sheet that was active/visible when the workbook was
# last saved — keep it that way.
$fa.Activate()
